{"js": "// Auto-generated: replace each math-problem cell's text (20 rows x 5 cols)\n// with the new value, in document order, matching the target diff.\nconst NEW_VALUES = [\"88-0=\", \"44-21=\", \"68-45=\", \"1+98=\", \"47-5=\", \"18+48=\", \"63+14=\", \"37-4=\", \"27+21=\", \"40-7=\", \"89-35=\", \"53+6=\", \"41+9=\", \"25+35=\", \"58+12=\", \"86+7=\", \"46-22=\", \"15+53=\", \"98-19=\", \"59+31=\", \"61-34=\", \"25+50=\", \"64+19=\", \"30+38=\", \"67+9=\", \"83-81=\", \"92-4=\", \"82+16=\", \"25+72=\", \"61+12=\", \"31+36=\", \"91-25=\", \"43-7=\", \"42-35=\", \"12+75=\", \"83-78=\", \"31+11=\", \"87-3=\", \"44-41=\", \"21+73=\", \"23+17=\", \"75-58=\", \"46+7=\", \"87-26=\", \"72-9=\", \"6+6=\", \"87-71=\", \"60-24=\", \"71-25=\", \"96-83=\", \"42-6=\", \"86-82=\", \"91-58=\", \"70+23=\", \"62-19=\", \"74-12=\", \"56-37=\", \"95-45=\", \"39+6=\", \"0+36=\", \"87-4=\", \"56-12=\", \"63+3=\", \"52+32=\", \"77-76=\", \"39+39=\", \"1+58=\", \"3+65=\", \"47-14=\", \"41-32=\", \"36+61=\", \"91+3=\", \"85-0=\", \"46+8=\", \"50+19=\", \"97-60=\", \"77+1=\", \"37-3=\", \"27+10=\", \"46+27=\", \"39+6=\", \"51-36=\", \"42+14=\", \"97-85=\", \"40-28=\", \"12+22=\", \"30-14=\", \"2+30=\", \"73-25=\", \"98-12=\", \"43+49=\", \"42+27=\", \"71-31=\", \"87-8=\", \"64-51=\", \"9+58=\", \"85-37=\", \"89-76=\", \"89-15=\", \"7+72=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst rowCount = rows.items.length;\n\n// Load cell counts per row first\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  const cells = rows.items[r].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    if (idx >= NEW_VALUES.length) break;\n    const cell = cells[c];\n    const cellRange = cell.body.getRange();\n    cellRange.insertText(NEW_VALUES[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated: replace each math-problem cell's text (20 rows x 5 cols)\n# with the new value, in document order, matching the target diff.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$values = @(\"88-0=\", \"44-21=\", \"68-45=\", \"1+98=\", \"47-5=\", \"18+48=\", \"63+14=\", \"37-4=\", \"27+21=\", \"40-7=\", \"89-35=\", \"53+6=\", \"41+9=\", \"25+35=\", \"58+12=\", \"86+7=\", \"46-22=\", \"15+53=\", \"98-19=\", \"59+31=\", \"61-34=\", \"25+50=\", \"64+19=\", \"30+38=\", \"67+9=\", \"83-81=\", \"92-4=\", \"82+16=\", \"25+72=\", \"61+12=\", \"31+36=\", \"91-25=\", \"43-7=\", \"42-35=\", \"12+75=\", \"83-78=\", \"31+11=\", \"87-3=\", \"44-41=\", \"21+73=\", \"23+17=\", \"75-58=\", \"46+7=\", \"87-26=\", \"72-9=\", \"6+6=\", \"87-71=\", \"60-24=\", \"71-25=\", \"96-83=\", \"42-6=\", \"86-82=\", \"91-58=\", \"70+23=\", \"62-19=\", \"74-12=\", \"56-37=\", \"95-45=\", \"39+6=\", \"0+36=\", \"87-4=\", \"56-12=\", \"63+3=\", \"52+32=\", \"77-76=\", \"39+39=\", \"1+58=\", \"3+65=\", \"47-14=\", \"41-32=\", \"36+61=\", \"91+3=\", \"85-0=\", \"46+8=\", \"50+19=\", \"97-60=\", \"77+1=\", \"37-3=\", \"27+10=\", \"46+27=\", \"39+6=\", \"51-36=\", \"42+14=\", \"97-85=\", \"40-28=\", \"12+22=\", \"30-14=\", \"2+30=\", \"73-25=\", \"98-12=\", \"43+49=\", \"42+27=\", \"71-31=\", \"87-8=\", \"64-51=\", \"9+58=\", \"85-37=\", \"89-76=\", \"89-15=\", \"7+72=\")\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$idx]\n        $idx++\n    }\n}\n"}
